# Weekly update: insert a new price record at the top of the Perejil /
# Terminal Hortofrutícola Agro Chillán block (row 107), pushing the
# existing rows 107-113 down to 108-114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 107 - this shifts rows
# 107:113 down to 108:114, carrying their values/styles with them.
$ws.Range("A107").EntireRow.Insert()

# Populate the newly inserted row 107 with this week's data.
$ws.Range("A107").Value = 7
$ws.Range("B107").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C107").Value = "Ñuble"
$ws.Range("D107").Value = 45166
$ws.Range("E107").Value = 16
$ws.Range("F107").Value = 100112044
$ws.Range("G107").Value = "Perejil"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 180
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1500
$ws.Range("N107").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O107").Value = "Región de Ñuble"
$ws.Range("P107").Value = 1500
$ws.Range("Q107").Value = 1
$ws.Range("R107").Value = "Hortaliza"
